$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.049038
$ws.Range("H2").Value = 63.147114
$ws.Range("I2").Value = 0.384846371905728
$ws.Range("J2").Value = 0.384846371905728
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 28.93895166666666
$ws.Range("N2").Value = 86.81685499999999
$ws.Range("O2").Value = 0.2389168411431201
$ws.Range("P2").Value = 0.2499063428956389
$ws.Range("Q2").Value = 609.13709331183
$ws.Range("R2").Value = 5482.233839806469
$ws.Range("S2").Value = 0.09194627950110693
$ws.Range("T2").Value = 0.09617554937961542
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.049038
$ws.Range("H3").Value = 63.147114
$ws.Range("I3").Value = 0.384846371905728
$ws.Range("J3").Value = 0.384846371905728
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.705903
$ws.Range("N3").Value = 107.117709
$ws.Range("O3").Value = 0.2947840562154431
$ws.Range("P3").Value = 0.3083432925040795
$ws.Range("Q3").Value = 751.574909071314
$ws.Range("R3").Value = 6764.174181641825
$ws.Range("S3").Value = 0.1134465745301674
$ws.Range("T3").Value = 0.1186647974216617
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.049038
$ws.Range("H4").Value = 63.147114
$ws.Range("I4").Value = 0.384846371905728
$ws.Range("J4").Value = 0.384846371905728
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 19.284198
$ws.Range("N4").Value = 57.852594
$ws.Range("O4").Value = 0.1592082437265831
$ws.Range("P4").Value = 0.1665313745074753
$ws.Range("Q4").Value = 405.913816501524
$ws.Range("R4").Value = 3653.224348513716
$ws.Range("S4").Value = 0.06127071497565838
$ws.Range("T4").Value = 0.06408899528767592
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 21.049038
$ws.Range("H5").Value = 63.147114
$ws.Range("I5").Value = 0.384846371905728
$ws.Range("J5").Value = 0.384846371905728
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.21726266666667
$ws.Range("N5").Value = 63.651788
$ws.Range("O5").Value = 0.1751674156138409
$ws.Range("P5").Value = 0.1832246233504832
$ws.Range("Q5").Value = 446.602968126648
$ws.Range("R5").Value = 4019.426713139832
$ws.Range("S5").Value = 0.06741254437508942
$ws.Range("T5").Value = 0.07051333154022699
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 21.049038
$ws.Range("H6").Value = 63.147114
$ws.Range("I6").Value = 0.384846371905728
$ws.Range("J6").Value = 0.384846371905728
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 15.9793095
$ws.Range("N6").Value = 31.958619
$ws.Range("O6").Value = 0.1319234433010128
$ws.Range("P6").Value = 0.09199436674232302
$ws.Range("Q6").Value = 336.349092879261
$ws.Range("R6").Value = 2018.094557275566
$ws.Range("S6").Value = 0.05077025852370579
$ws.Range("T6").Value = 0.03540369827654798
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 33.645613
$ws.Range("H7").Value = 100.936839
$ws.Range("I7").Value = 0.615153628094272
$ws.Range("J7").Value = 0.615153628094272
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 28.93895166666666
$ws.Range("N7").Value = 86.81685499999999
$ws.Range("O7").Value = 0.2389168411431201
$ws.Range("P7").Value = 0.2499063428956389
$ws.Range("Q7").Value = 973.6687684023717
$ws.Range("R7").Value = 8763.018915621344
$ws.Range("S7").Value = 0.1469705616420132
$ws.Range("T7").Value = 0.1537307935160235
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 33.645613
$ws.Range("H8").Value = 100.936839
$ws.Range("I8").Value = 0.615153628094272
$ws.Range("J8").Value = 0.615153628094272
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 35.705903
$ws.Range("N8").Value = 107.117709
$ws.Range("O8").Value = 0.2947840562154431
$ws.Range("P8").Value = 0.3083432925040795
$ws.Range("Q8").Value = 1201.346994153539
$ws.Range("R8").Value = 10812.12294738185
$ws.Range("S8").Value = 0.1813374816852756
$ws.Range("T8").Value = 0.1896784950824179
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 33.645613
$ws.Range("H9").Value = 100.936839
$ws.Range("I9").Value = 0.615153628094272
$ws.Range("J9").Value = 0.615153628094272
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.284198
$ws.Range("N9").Value = 57.852594
$ws.Range("O9").Value = 0.1592082437265831
$ws.Range("P9").Value = 0.1665313745074753
$ws.Range("Q9").Value = 648.8286629233741
$ws.Range("R9").Value = 5839.457966310366
$ws.Range("S9").Value = 0.09793752875092471
$ws.Range("T9").Value = 0.1024423792197994
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 33.645613
$ws.Range("H10").Value = 100.936839
$ws.Range("I10").Value = 0.615153628094272
$ws.Range("J10").Value = 0.615153628094272
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 21.21726266666667
$ws.Range("N10").Value = 63.651788
$ws.Range("O10").Value = 0.1751674156138409
$ws.Range("P10").Value = 0.1832246233504832
$ws.Range("Q10").Value = 713.8678086020147
$ws.Range("R10").Value = 6424.810277418132
$ws.Range("S10").Value = 0.1077548712387515
$ws.Range("T10").Value = 0.1127112918102562
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 33.645613
$ws.Range("H11").Value = 100.936839
$ws.Range("I11").Value = 0.615153628094272
$ws.Range("J11").Value = 0.615153628094272
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 15.9793095
$ws.Range("N11").Value = 31.958619
$ws.Range("O11").Value = 0.1319234433010128
$ws.Range("P11").Value = 0.09199436674232302
$ws.Range("Q11").Value = 537.6336634442235
$ws.Range("R11").Value = 3225.801980665341
$ws.Range("S11").Value = 0.081153184777307
$ws.Range("T11").Value = 0.05659066846577503
